$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Fix header typo: "Author" -> "Authors" on the "R2T.csv" sheet (A1), which
#    is the sheet whose row 2 actually lists multiple authors.
# ---------------------------------------------------------------------------
$wsR2T = $wb.Worksheets.Item("R2T.csv")
$wsR2T.Range("A1").Value = "Authors"

# ---------------------------------------------------------------------------
# 2) Add a new worksheet "TRm.csv" as the last sheet, a test case with
#    pre-existing data in the abstract/doi columns.
#    We copy the last existing sheet ("TR1DA.csv(2)") so that all of the
#    page setup / formatting boilerplate (sheetPr, cols, sheetProtection,
#    printOptions, pageMargins, pageSetup, headerFooter, Print_Area name,
#    etc.) is carried over automatically, matching the rest of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "TRm.csv"

# Row 2: pre-existing (already filled-in) abstract/doi values.
$newSheet.Range("C2").Value = "Not Sure"
$newSheet.Range("D2").Value = "Couldn't find it"

# Row 3: same title/author as row 2, but with no DOI/Abstract filled in yet.
$newSheet.Range("A3").Value = $newSheet.Range("A2").Value2
$newSheet.Range("B3").Value = $newSheet.Range("B2").Value2

# The copy operation already duplicated the "_xlnm.Print_Area" defined name
# for this sheet; also add the (gnumeric-specific) "_xlnm.Sheet_Title" name
# that every other sheet has.
$newSheet.Names.Add("_xlnm.Sheet_Title", "=""TRm.csv""") | Out-Null

# Restore the originally active sheet/tab (unchanged by this edit).
$wb.Worksheets.Item(1).Activate()
